$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("C#")

# --- New "STRING specifics" notes block (rows 93-128) ---
# Note: row 98's A/B/C cells are written with A98 last (matches the
# original authoring order: the method-name label was filled in after
# the description/sample text that already used that exact phrase, so
# the shared-string table gains "concatenates strings method" before
# "string.Contact()").
$ws.Range("A93").Value = "STRING specifics"
$ws.Range("A95").Value = "char[] charArray = var.ToCharArray()"
$ws.Range("B95").Value = "converting string into char array"
$ws.Range("C95").Value = "string text = `"text`""
$ws.Range("C96").Value = "char[] charArray = text.ToCharArray() // ['t', 'e', 'x', 't']"
$ws.Range("B98").Value = "concatenates strings method"
$ws.Range("C98").Value = "string greet = `"hello, `"    string name = `"George`""
$ws.Range("C99").Value = "string result = string.Contact( greet, name)  // Hello, George"
$ws.Range("A98").Value = "string.Contact()"
$ws.Range("A101").Value = "IndexOf()"
$ws.Range("B101").Value = "returns the first match index or -1"
$ws.Range("C101").Value = "string fruits = `"banana, apple, kiwi`""
$ws.Range("C102").Value = "(fruits.IndexOf(`"banana`")) // 0    first letter b is on index# 0"
$ws.Range("A104").Value = "LastIndexOf()"
$ws.Range("B104").Value = "returns last match index occurance"
$ws.Range("C104").Value = "string fruits = `"banana, apple, kiwi, banana, apple`""
$ws.Range("C105").Value = "(fruits.IndexOf(`"banana`")) //21    first letter b is on index# 0 from back to forward"
$ws.Range("A107").Value = "Contains()"
$ws.Range("B107").Value = "bool. Finds a string if is contained in a string var"
$ws.Range("C107").Value = "string text = `"I love fruits`""
$ws.Range("C108").Value = "(text.Contains(`"fruits`"))  - will return True as the string contains the word"
$ws.Range("A110").Value = "Substring(int startIndex, int length)"
$ws.Range("C110").Value = "string card = `"10C`""
$ws.Range("C111").Value = "string power = card.Substring(0, 2)  // 10"
$ws.Range("A113").Value = "Substring(int startIndex)"
$ws.Range("C113").Value = "string name = `"My name is John`"   (name starts at the 11th index)"
$ws.Range("C114").Value = "string exactName = string.Substring(11)  // returns John"
$ws.Range("A116").Value = "Replace(match, replacement)"
$ws.Range("B116").Value = "replaces a matched string with a new one"
$ws.Range("C116").Value = ".Replace(ggajev1@abv.bg, gadzhev@abv.bg) // new string will be gadzhev@abv.bg"
$ws.Range("C117").Value = "replacement text should be same string length"
$ws.Range("A120").Value = "StringBuilder sb = new StringBuilder()"
$ws.Range("B120").Value = "build/modify strings class"
$ws.Range("A122").Value = "sb.Append"
$ws.Range("A126").Value = "Stopwatch sw = new Stopwatch()"
$ws.Range("A127").Value = "sw.Start()"
$ws.Range("A128").Value = "(sw.ElapsedMilliSeconds)"

# Column C got wider to fit the longer code samples/notes added above.
$ws.Columns.Item(3).ColumnWidth = 76.6

# Built-in "Normal" cell style was re-localized from Bulgarian to English.
$normalStyle = $wb.Styles.Item(1)
$normalStyle.Name = "Normal"

# Leave the selection where the author left off editing (matches the
# saved cursor position after typing the last note), and scroll the
# window down so row 97 is pinned at the top.
$excel.ActiveWindow.ScrollRow = 97
$ws.Range("B127").Select() | Out-Null
